$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L2: numeric 10 -> 80
$ws.Range("L2").Value = 80

# L3: numeric 0 -> text "20" (stored as text, not a number, no residual
# number-format styling left on the cell)
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "20"
$ws.Range("L3").ClearFormats()
